# Apply the "sua luong tai tu vsign sua lai drive" changes to the Config sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Config")

# 1) Update "stage" value (row 7, column B): 1 -> 2
$ws.Range("B7").Value = 2

# 2) Update "createdDay" value (row 8, column B): 07/28/2023 12:20:09 -> 08/03/2023 10:14:47
$ws.Range("B8").Value = "08/03/2023 10:14:47"

# 3) Update "uploadFolderName" value (row 9, column B): 28-07-2023 -> DanhSachHHVT 03/08/2023
$ws.Range("B9").Value = "DanhSachHHVT 03/08/2023"

# 4) Remove the "nccFileAddress" row entirely (old row 14), shifting everything below it up by one.
$ws.Rows("14:14").Delete()

# 5) Reset the view: top-left cell back to A1 and selection on B7.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B7").Select()

$wb.Save()
